$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 1).Value = -0.10522383444921957
$ws.Cells.Item(1, 2).Value = 0.10479246395502173
$ws.Cells.Item(2, 1).Value = -0.015034255586574119
$ws.Cells.Item(2, 2).Value = 0.014065861576341376
$ws.Cells.Item(3, 1).Value = 0.08886872569704352
$ws.Cells.Item(3, 2).Value = -0.08941723635071597
$ws.Cells.Item(4, 1).Value = -0.19457497554336456
$ws.Cells.Item(4, 2).Value = 0.19333456897670587
$ws.Cells.Item(5, 1).Value = -0.18733456952467886
$ws.Cells.Item(5, 2).Value = 0.18482026627278625
$ws.Cells.Item(6, 1).Value = -0.059927525861781294
$ws.Cells.Item(6, 2).Value = 0.05987733832015696
$ws.Cells.Item(7, 1).Value = -0.0398773389771403
$ws.Cells.Item(7, 2).Value = 0.03979687056066972
$ws.Cells.Item(8, 1).Value = -0.019796871220780332
$ws.Cells.Item(8, 2).Value = 0.01977014104207342
$ws.Cells.Item(9, 1).Value = -0.013770141617934328
$ws.Cells.Item(9, 2).Value = 0.013757291101064162
$ws.Cells.Item(10, 1).Value = -0.007757291678672118
$ws.Cells.Item(10, 2).Value = 0.007758330278413439
$ws.Cells.Item(11, 1).Value = -0.00325833084681193
$ws.Cells.Item(11, 2).Value = 0.0032587399175980636
$ws.Cells.Item(12, 1).Value = -0.049121472950132805
$ws.Cells.Item(12, 2).Value = 0.04866809728676014
$ws.Cells.Item(13, 1).Value = -0.03915133565518758
$ws.Cells.Item(13, 2).Value = 0.03908484909886489
$ws.Cells.Item(14, 1).Value = -0.027084849726634275
$ws.Cells.Item(14, 2).Value = 0.02705280235262375
$ws.Cells.Item(15, 1).Value = -0.021052802944928395
$ws.Cells.Item(15, 2).Value = 0.021027581245259164
$ws.Cells.Item(16, 1).Value = -0.015027581839187398
$ws.Cells.Item(16, 2).Value = 0.015004403272817601
$ws.Cells.Item(17, 1).Value = -0.00900440386904311
$ws.Cells.Item(17, 2).Value = 0.008999999384767143
$ws.Cells.Item(18, 1).Value = -0.036109943106982456
$ws.Cells.Item(18, 2).Value = 0.03609664379099087
$ws.Cells.Item(19, 1).Value = -0.027096644350875998
$ws.Cells.Item(19, 2).Value = 0.02701366810978234
$ws.Cells.Item(20, 1).Value = -0.018013668674269567
$ws.Cells.Item(20, 2).Value = 0.01800427511715874
$ws.Cells.Item(21, 1).Value = -0.009004275682231722
$ws.Cells.Item(21, 2).Value = 0.008999999434542438
$ws.Cells.Item(22, 1).Value = -0.09393458783425146
$ws.Cells.Item(22, 2).Value = 0.093625899650565
$ws.Cells.Item(23, 1).Value = -0.08462590021492211
$ws.Cells.Item(23, 2).Value = 0.0841250609509725
$ws.Cells.Item(24, 1).Value = -0.042125061728516044
$ws.Cells.Item(24, 2).Value = 0.04199999921859554
$ws.Cells.Item(25, 1).Value = -0.09489754001896245
$ws.Cells.Item(25, 2).Value = 0.09465601711839255
$ws.Cells.Item(26, 1).Value = -0.08865601768914289
$ws.Cells.Item(26, 2).Value = 0.08834436764047382
$ws.Cells.Item(27, 1).Value = -0.08234436821350011
$ws.Cells.Item(27, 2).Value = 0.08127646004388467
$ws.Cells.Item(28, 1).Value = -0.07527646062873217
$ws.Cells.Item(28, 2).Value = 0.07453654229329931
$ws.Cells.Item(29, 1).Value = -0.06253654292296851
$ws.Cells.Item(29, 2).Value = 0.062169901445164655
$ws.Cells.Item(30, 1).Value = -0.04216990212840921
$ws.Cells.Item(30, 2).Value = 0.04201956123942585
$ws.Cells.Item(31, 1).Value = -0.027019561897143163
$ws.Cells.Item(31, 2).Value = 0.02700063533537289
$ws.Cells.Item(32, 1).Value = -0.006000636031154194
$ws.Cells.Item(32, 2).Value = 0.005999999397076294
